# Updates match-result / odds data for the "Germany Landesliga" sheet.
# Source data refresh caused a handful of fixture rows to be re-ordered;
# for each affected pair of rows the match id (column B) together with
# the away team (column G) and every score/odds column (H through AC)
# move from one row to the other, while the row's sequence id (column A)
# and home team (column F) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2) {
    # Column B (match id)
    $cols = @(2) + (7..29)   # B, then G..AC

    foreach ($col in $cols) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)

        $val1 = $cell1.Value2
        $val2 = $cell2.Value2

        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}

Swap-RowData 16 17
Swap-RowData 29 30
Swap-RowData 37 38
Swap-RowData 48 49
